$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix bug in driver: rows 6 and 7 previously held placeholder "test" data.
# Row 6 -> Steak, Row 7 -> Noodles, with corrected price/special-item values.
$ws.Range("B6").Value = "Steak"
$ws.Range("C6").Value = 5.0
$ws.Range("D6").Value = $true

$ws.Range("B7").Value = "Noodles"
$ws.Range("C7").Value = 3.0

# Row 8 (ice cream) is now marked Active.
$ws.Range("E8").Value = $true

# Complete the chef page: add two new menu items.
$ws.Range("A9").Value = 8.0
$ws.Range("B9").Value = "Naan"
$ws.Range("C9").Value = 2.0
$ws.Range("D9").Value = $false
$ws.Range("E9").Value = $true

$ws.Range("A10").Value = 9.0
$ws.Range("B10").Value = "Biriyani"
$ws.Range("C10").Value = 4.0
$ws.Range("D10").Value = $false
$ws.Range("E10").Value = $true
